# Apply the restructuring of columns I..P described by the commit diff.
#
# Header row (row 1): rename columns I..P to their new labels.
# Data rows (2..65): the old "Durée*" counter columns (I,J,K) are dropped,
# and the old Début/Elaboration/CTCQ/Conformité/Approbation values (L,M,N,O,P)
# slide left into the new layout:
#   new I (Début)                      <- old L (date)
#   new J (Elaboration Prévisionnelle) <- old M (date)
#   new K (Elaboration Effective)      <- (new, empty)
#   new L (CTCQ Prévisionnelle)        <- old N (text)
#   new M (CTCQ Effective)             <- (new, empty)
#   new N (Conformité)                 <- old O (text)
#   new O (Approbation Prévisionnelle) <- old P (text)
#   new P (Approbation Effective)      <- (new, empty)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row ----
$ws.Range("I1").Value = "Début"
$ws.Range("J1").Value = "Elaboration Prévisionnelle"
$ws.Range("K1").Value = "Elaboration Effective"
$ws.Range("L1").Value = "CTCQ Prévisionnelle"
$ws.Range("M1").Value = "CTCQ Effective"
$ws.Range("N1").Value = "Conformité"
$ws.Range("O1").Value = "Approbation Prévisionnelle"
$ws.Range("P1").Value = "Approbation Effective"

# ---- Data rows ----
$lastRow = $ws.Range("A1").End(4).Row  # xlDown = 4, find the bottom of column A
if ($lastRow -lt 65) { $lastRow = 65 }

for ($r = 2; $r -le $lastRow; $r++) {

    $cellL = $ws.Cells.Item($r, 12)   # L: old Début (date)
    $cellM = $ws.Cells.Item($r, 13)   # M: old Elaboration (date)
    $cellN = $ws.Cells.Item($r, 14)   # N: old CTCQ (text)
    $cellO = $ws.Cells.Item($r, 15)   # O: old Conformité (text)
    $cellP = $ws.Cells.Item($r, 16)   # P: old Approbation (text)

    # Capture old values/formats before anything in this row gets overwritten.
    $lVal = $cellL.Value2
    $lFmt = $cellL.NumberFormat
    $mVal = $cellM.Value2
    $mFmt = $cellM.NumberFormat

    $nText = $cellN.Text
    $oText = $cellO.Text
    $pText = $cellP.Text

    $cellI = $ws.Cells.Item($r, 9)
    $cellJ = $ws.Cells.Item($r, 10)
    $cellK = $ws.Cells.Item($r, 11)
    $cellNewL = $ws.Cells.Item($r, 12)
    $cellNewM = $ws.Cells.Item($r, 13)
    $cellNewN = $ws.Cells.Item($r, 14)
    $cellNewO = $ws.Cells.Item($r, 15)
    $cellNewP = $ws.Cells.Item($r, 16)

    # new I <- old L (date, keep the numeric date format)
    $cellI.NumberFormat = $lFmt
    $cellI.Value = $lVal

    # new J <- old M (date, keep the numeric date format)
    $cellJ.NumberFormat = $mFmt
    $cellJ.Value = $mVal

    # new K <- empty (Elaboration Effective is not yet populated)
    $cellK.ClearFormats()
    $cellK.ClearContents()

    # new L <- old N, as plain text (leading ' keeps numeric-looking text as text)
    $cellNewL.ClearFormats()
    if ($nText -eq "") {
        $cellNewL.ClearContents()
    } else {
        $cellNewL.Value = "'" + $nText
    }

    # new M <- empty (CTCQ Effective is not yet populated)
    $cellNewM.ClearFormats()
    $cellNewM.ClearContents()

    # new N <- old O, as plain text
    $cellNewN.ClearFormats()
    if ($oText -eq "") {
        $cellNewN.ClearContents()
    } else {
        $cellNewN.Value = "'" + $oText
    }

    # new O <- old P, as plain text
    $cellNewO.ClearFormats()
    if ($pText -eq "") {
        $cellNewO.ClearContents()
    } else {
        $cellNewO.Value = "'" + $pText
    }

    # new P <- empty (Approbation Effective is not yet populated)
    $cellNewP.ClearFormats()
    $cellNewP.ClearContents()
}
